$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the old last row (42) first, shifting rows up and updating dimension
$ws.Range("A42:K42").Delete() | Out-Null

$data = New-Object 'object[,]' 27,11
$data[0,0] = 'Gaur et al. (2022)'
$data[0,1] = 'IE'
$data[0,2] = 'ILED'
$data[0,3] = 'Freight'
$data[0,4] = 'transported goods per capita and year'
$data[0,5] = 2050
$data[0,6] = 'tkm/cap/year'
$data[0,7] = 1525
$data[0,8] = 2295
$data[0,9] = 0.6644880174291939
$data[0,10] = 1
$data[1,0] = 'négaWatt et al. (2023)'
$data[1,1] = 'EU27'
$data[1,2] = 'EU path'
$data[1,3] = 'Transport'
$data[1,4] = 'transported goods per capita and year'
$data[1,5] = 2050
$data[1,6] = 'tkm/cap/year'
$data[1,7] = 5291.85
$data[1,8] = 2295
$data[1,9] = 2.305816993464052
$data[1,10] = 561
$data[2,0] = 'RTE (2021)'
$data[2,1] = 'FR'
$data[2,2] = 'Sufficiency'
$data[2,3] = 'Transport'
$data[2,4] = 'transported goods per capita and year'
$data[2,5] = 2050
$data[2,6] = 'tkm/cap/year'
$data[2,7] = 5434.78
$data[2,8] = 2295
$data[2,9] = 2.368095860566449
$data[2,10] = 283
$data[3,0] = 'UBA (2020)'
$data[3,1] = 'DE'
$data[3,2] = 'GreenSupreme'
$data[3,3] = 'Freight'
$data[3,4] = 'transported goods per capita and year'
$data[3,5] = 2050
$data[3,6] = 'tkm/cap/year'
$data[3,7] = 8122.17
$data[3,8] = 2295
$data[3,9] = 3.539071895424837
$data[3,10] = 132
$data[4,0] = 'ADEME (2022)'
$data[4,1] = 'FR'
$data[4,2] = 'S1: Low production'
$data[4,3] = 'Transport'
$data[4,4] = 'average distance travelled per capita and year'
$data[4,5] = 2050
$data[4,6] = 'pkm/cap/year'
$data[4,7] = 10672
$data[4,8] = 10000
$data[4,9] = 1.0672
$data[4,10] = 148
$data[5,0] = 'ADEME (2022)'
$data[5,1] = 'FR'
$data[5,2] = 'S1: Low production'
$data[5,3] = 'Transport'
$data[5,4] = 'average distance travelled per capita and year'
$data[5,5] = 2050
$data[5,6] = 'pkm/cap/year'
$data[5,7] = 12081.5
$data[5,8] = 10000
$data[5,9] = 1.20815
$data[5,10] = 149
$data[6,0] = 'CTI 2050 Roadmap Tool (2018)'
$data[6,1] = 'EU28'
$data[6,2] = 'demand'
$data[6,3] = 'Transport'
$data[6,4] = 'average distance travelled per capita and year'
$data[6,5] = 2050
$data[6,6] = 'pkm/cap/year'
$data[6,7] = 10450.25
$data[6,8] = 10000
$data[6,9] = 1.045025
$data[6,10] = 387
$data[7,0] = 'Gaur et al. (2022)'
$data[7,1] = 'IE'
$data[7,2] = 'ILED'
$data[7,3] = 'Mobility'
$data[7,4] = 'average distance travelled per capita and year'
$data[7,5] = 2050
$data[7,6] = 'pkm/cap/year'
$data[7,7] = 12000
$data[7,8] = 10000
$data[7,9] = 1.2
$data[7,10] = 2
$data[8,0] = 'négaWatt (2022)'
$data[8,1] = 'FR'
$data[8,2] = '2050'
$data[8,3] = 'Mobility'
$data[8,4] = 'average distance travelled per capita and year'
$data[8,5] = 2050
$data[8,6] = 'pkm/cap/year'
$data[8,7] = 14500
$data[8,8] = 10000
$data[8,9] = 1.45
$data[8,10] = 25
$data[9,0] = 'négaWatt et al. (2023)'
$data[9,1] = 'EU27'
$data[9,2] = 'EU path'
$data[9,3] = 'Transport'
$data[9,4] = 'average distance travelled per capita and year'
$data[9,5] = 2050
$data[9,6] = 'pkm/cap/year'
$data[9,7] = 13803.44
$data[9,8] = 10000
$data[9,9] = 1.380344
$data[9,10] = 554
$data[10,0] = 'RTE (2021)'
$data[10,1] = 'FR'
$data[10,2] = 'Sufficiency'
$data[10,3] = 'Transport'
$data[10,4] = 'average distance travelled per capita and year'
$data[10,5] = 2050
$data[10,6] = 'pkm/cap/year'
$data[10,7] = 13550
$data[10,8] = 10000
$data[10,9] = 1.355
$data[10,10] = 281
$data[11,0] = 'UBA (2020)'
$data[11,1] = 'DE'
$data[11,2] = 'GreenSupreme'
$data[11,3] = 'Mobility'
$data[11,4] = 'average distance travelled per capita and year'
$data[11,5] = 2050
$data[11,6] = 'pkm/cap/year'
$data[11,7] = 13337.6
$data[11,8] = 10000
$data[11,9] = 1.33376
$data[11,10] = 91
$data[12,0] = 'UBA (2020)'
$data[12,1] = 'DE'
$data[12,2] = 'GreenSupreme'
$data[12,3] = 'Mobility'
$data[12,4] = 'average distance travelled per capita and year'
$data[12,5] = 2050
$data[12,6] = 'pkm/cap/year'
$data[12,7] = 16019.03
$data[12,8] = 10000
$data[12,9] = 1.601903
$data[12,10] = 92
$data[13,0] = 'van de Ven et al. (2018)'
$data[13,1] = 'EU27'
$data[13,2] = 'Enthusiastic profile'
$data[13,3] = 'Mobility'
$data[13,4] = 'average distance travelled per capita and year'
$data[13,5] = 2050
$data[13,6] = 'pkm/cap/year'
$data[13,7] = 11829.32
$data[13,8] = 10000
$data[13,9] = 1.182932
$data[13,10] = 366
$data[14,0] = 'ADEME (2022)'
$data[14,1] = 'FR'
$data[14,2] = 'S1: Low production'
$data[14,3] = 'Buildings'
$data[14,4] = 'per capita floor area in commercial and public buildings'
$data[14,5] = 2050
$data[14,6] = 'm2/cap'
$data[14,7] = 12
$data[14,8] = 5.4
$data[14,9] = 2.222222222222222
$data[14,10] = 160
$data[15,0] = 'Gaur et al. (2022)'
$data[15,1] = 'IE'
$data[15,2] = 'ILED'
$data[15,3] = 'Buildings'
$data[15,4] = 'per capita floor area in commercial and public buildings'
$data[15,5] = 2050
$data[15,6] = 'm2/cap'
$data[15,7] = 16
$data[15,8] = 5.4
$data[15,9] = 2.962962962962963
$data[15,10] = 22
$data[16,0] = 'négaWatt (2022)'
$data[16,1] = 'FR'
$data[16,2] = '2050'
$data[16,3] = 'Buildings'
$data[16,4] = 'per capita floor area in commercial and public buildings'
$data[16,5] = 2050
$data[16,6] = 'm2/cap'
$data[16,7] = 15.14
$data[16,8] = 5.4
$data[16,9] = 2.803703703703704
$data[16,10] = 37
$data[17,0] = 'RTE (2021)'
$data[17,1] = 'FR'
$data[17,2] = 'Sufficiency'
$data[17,3] = 'Buildings'
$data[17,4] = 'per capita floor area in commercial and public buildings'
$data[17,5] = 2050
$data[17,6] = 'm2/capita'
$data[17,7] = 16.48
$data[17,8] = 5.4
$data[17,9] = 3.051851851851852
$data[17,10] = 272
$data[18,0] = 'ADEME (2022)'
$data[18,1] = 'FR'
$data[18,2] = 'S1: Low production'
$data[18,3] = 'Buildings'
$data[18,4] = 'living space per capita'
$data[18,5] = 2050
$data[18,6] = 'm2/cap'
$data[18,7] = 48
$data[18,8] = 15
$data[18,9] = 3.2
$data[18,10] = 178
$data[19,0] = 'Eerma et al. (2022)'
$data[19,1] = 'DE'
$data[19,2] = 'High Ambition'
$data[19,3] = 'Buildings'
$data[19,4] = 'living space per capita'
$data[19,5] = 2050
$data[19,6] = 'm^2/person'
$data[19,7] = 30
$data[19,8] = 15
$data[19,9] = 2
$data[19,10] = 196
$data[20,0] = 'Fishman et al. (2021)'
$data[20,1] = 'DE'
$data[20,2] = 'based on LED'
$data[20,3] = ""
$data[20,4] = 'living space per capita'
$data[20,5] = 2050
$data[20,6] = 'm2/cap'
$data[20,7] = 32
$data[20,8] = 15
$data[20,9] = 2.133333333333333
$data[20,10] = 83
$data[21,0] = 'Gaur et al. (2022)'
$data[21,1] = 'IE'
$data[21,2] = 'ILED'
$data[21,3] = 'Buildings'
$data[21,4] = 'living space per capita'
$data[21,5] = 2050
$data[21,6] = 'm2/cap'
$data[21,7] = 43.43
$data[21,8] = 15
$data[21,9] = 2.895333333333333
$data[21,10] = 11
$data[22,0] = 'négaWatt (2022)'
$data[22,1] = 'FR'
$data[22,2] = '2050'
$data[22,3] = 'Buildings'
$data[22,4] = 'living space per capita'
$data[22,5] = 2050
$data[22,6] = 'm2/cap'
$data[22,7] = 39.47
$data[22,8] = 15
$data[22,9] = 2.631333333333333
$data[22,10] = 35
$data[23,0] = 'négaWatt et al. (2023)'
$data[23,1] = 'EU27'
$data[23,2] = 'EU path'
$data[23,3] = 'Buildings'
$data[23,4] = 'living space per capita'
$data[23,5] = 2050
$data[23,6] = 'm²/cap'
$data[23,7] = 41.49
$data[23,8] = 15
$data[23,9] = 2.766
$data[23,10] = 552
$data[24,0] = 'UBA (2020)'
$data[24,1] = 'DE'
$data[24,2] = 'GreenSupreme'
$data[24,3] = 'Buildings'
$data[24,4] = 'living space per capita'
$data[24,5] = 2050
$data[24,6] = 'm2/cap'
$data[24,7] = 41.17
$data[24,8] = 15
$data[24,9] = 2.744666666666667
$data[24,10] = 114
$data[25,0] = 'van Sluisveld et al. (2020)'
$data[25,1] = 'EU'
$data[25,2] = 'RegChange'
$data[25,3] = 'Living space'
$data[25,4] = 'living space per capita'
$data[25,5] = 2050
$data[25,6] = 'm2/cap'
$data[25,7] = 43.3
$data[25,8] = 15
$data[25,9] = 2.886666666666666
$data[25,10] = 730
$data[26,0] = 'UBA (2020)'
$data[26,1] = 'DE'
$data[26,2] = 'GreenSupreme'
$data[26,3] = 'Agriculture'
$data[26,4] = 'meat consumption per capita and day'
$data[26,5] = 2050
$data[26,6] = 'g/cap/day'
$data[26,7] = 42.86
$data[26,8] = 62.25
$data[26,9] = 0.6885140562248996
$data[26,10] = 127

$ws.Range("A15:K41").Value = $data
